# Clean-up of input tables
# - Rename the worksheet from "updated" to "Tabelle1"
# - Update the active cell selection to AJ9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Tabelle1"

$ws.Range("AJ9").Select()
